# Append the new match row (row 80) to the HNL 2023-2024 sheet,
# mirroring the formatting of the previous data row (row 79).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 79
$row = 80

# Column A: running index - copy formatting (bold, border, centered) from the row above.
$ws.Cells.Item($srcRow, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item($row, 1).Value = 79

$ws.Cells.Item($row, 2).Value = "croatia"
$ws.Cells.Item($row, 3).Value = "hnl"
$ws.Cells.Item($row, 4).Value = "2023-2024"

# Column E: match date/time - copy the date/time number format from the row above.
$ws.Cells.Item($srcRow, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item($row, 5).Value = 45262.625

$ws.Cells.Item($row, 6).Value = "Osijek"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Rijeka"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 2.61
$ws.Cells.Item($row, 11).Value = "25/11/2023 18:12"
$ws.Cells.Item($row, 12).Value = 3.33
$ws.Cells.Item($row, 13).Value = "02/12/2023 13:57"
$ws.Cells.Item($row, 14).Value = 3.27
$ws.Cells.Item($row, 15).Value = "25/11/2023 18:12"
$ws.Cells.Item($row, 16).Value = 3.34
$ws.Cells.Item($row, 17).Value = "02/12/2023 13:57"
$ws.Cells.Item($row, 18).Value = 2.73
$ws.Cells.Item($row, 19).Value = "25/11/2023 18:12"
$ws.Cells.Item($row, 20).Value = 2.25
$ws.Cells.Item($row, 21).Value = "02/12/2023 12:59"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/osijek-rijeka/ptWZnCPc/"
